$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "68.791.27"
$ws.Range("E2").Value = "  +1.69%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.737.20"
$ws.Range("E3").Value = "  -1.85%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.04%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'602.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.94%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'167.59"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.00%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "3.736.92"
$ws.Range("E7").Value = "  -1.81%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.09%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  +2.00%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +2.39%  "

# Row 11 - Toncoin
$ws.Range("D11").Value = "'6.36"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.80%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  -0.92%  "

# Row 13 - Avalanche
$ws.Range("D13").Value = "'38.02"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.65%  "

# Row 14 - ShibaInu
$ws.Range("E14").Value = "  -0.22%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.361.95"
$ws.Range("E15").Value = "  -1.90%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "3.740.85"
$ws.Range("E16").Value = "  -1.96%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "68.757.62"
$ws.Range("E17").Value = "  +1.61%  "

# Row 18 - Polkadot
$ws.Range("E18").Value = "  +0.78%  "

# Row 19 - TRON
$ws.Range("E19").Value = "  +0.51%  "

# Row 20 - Chainlink
$ws.Range("D20").Value = "'17.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.69%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "'496.92"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.18%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "'10.11"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +11.50%  "

# Row 23 - Polygon
$ws.Range("D23").Value = "'0.723"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.92%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "'84.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.10%  "

# Row 25 - Fetch.AI
$ws.Range("E25").Value = "  -2.63%  "

# Row 26 - PEPE
$ws.Range("E26").Value = "  -7.73%  "

# Row 27 - InternetComputer(DFINITY)
$ws.Range("E27").Value = "  +0.66%  "

# Row 28 - RenderToken
$ws.Range("D28").Value = "'10.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.76%  "

# Row 29 - Dai
$ws.Range("E29").Value = "  +0.01%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  -0.20%  "

# Row 31 - ImmutableX
$ws.Range("E31").Value = "  +0.44%  "

# Row 32 - NEARProtocol
$ws.Range("E32").Value = "  +3.70%  "

# Row 33 - EthereumClassic
$ws.Range("D33").Value = "'31.68"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.62%  "

# Row 34 - WrappedeETH
$ws.Range("D34").Value = "3.882.46"
$ws.Range("E34").Value = "  -1.84%  "

# Row 35 - RenzoRestakedETH->Hedera
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "'0.108"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.25%  "

# Row 36 - Hedera->RenzoRestakedETH
$ws.Range("B36").Value = "RenzoRestakedETH"
$ws.Range("C36").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D36").Value = "3.669.48"
$ws.Range("E36").Value = "  -2.04%  "

# Row 37 - FirstDigitalUSD
$ws.Range("D37").Value = "'1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.00%  "

# Row 38 - Mantle
$ws.Range("D38").Value = "'1.02"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.98%  "

# Row 39 - Filecoin
$ws.Range("D39").Value = "'5.82"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.63%  "

# Row 40 - Kaspa
$ws.Range("E40").Value = "  -2.00%  "

# Row 41 - TheGraph
$ws.Range("E41").Value = "  -1.15%  "

# Row 42 - Bittensor
$ws.Range("D42").Value = "'433.76"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.48%  "

# Row 43 - OKB
$ws.Range("D43").Value = "'49.02"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.35%  "

# Row 44 - Stacks
$ws.Range("E44").Value = "  -1.10%  "

# Row 45 - dogwifhat
$ws.Range("E45").Value = "  -0.02%  "

# Row 46 - Cosmos
$ws.Range("D46").Value = "'8.42"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.04%  "

# Row 48 - Arweave
$ws.Range("D48").Value = "'40.54"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.47%  "

# Row 49 - Monero
$ws.Range("D49").Value = "'141.66"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.35%  "

# Row 50 - VeChain
$ws.Range("E50").Value = "  +1.18%  "

# Row 51 - Maker
$ws.Range("D51").Value = "2.744.11"
$ws.Range("E51").Value = "  -3.12%  "
